$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.956.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.333.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.29%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.27"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.25"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0797"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.90"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.337.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.811"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.887.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.76"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.00%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0912"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.83"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.45"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.89"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.54%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.53"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.76%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.10"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.63"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.80%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.59"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0725"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.85"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.80%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.88%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.11%  "

$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +14.75%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.004.06"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.99%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0284"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.51%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.16"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.63%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.59%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.73"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.559.74"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.18%  "
